# Update the "想去人数" (want-to-go count) values in column F
# for the rows that changed between the two scrapes.
#
# The same update applies identically to both the "展览" sheet
# and the "全部类型" sheet, since they contain the same data.

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F
$updates = @{
    3  = 96
    4  = 279
    6  = 574
    8  = 2037
    11 = 4419
    12 = 34
    14 = 102
    16 = 119
    17 = 29
    18 = 18
    19 = 75
    20 = 3247
    22 = 486
    23 = 18
    25 = 77
    26 = 89
    29 = 58
    30 = 201
    32 = 591
    33 = 1866
    34 = 289
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
